$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 06:38 PM"

# --- Top Gainers sheet: refresh leaderboard (drop oldest row, shift up, append new entry) ---
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Rows.Item(61).Delete()
$gainers.Cells.Item(61, 1).Value = "🚀"
$gainers.Cells.Item(61, 2).Value = "APARINDS"
$gainers.Cells.Item(61, 3).Value = 3.8924
$gainers.Cells.Item(61, 4).Value = 8.3414
$gainers.Cells.Item(61, 5).Value = 15.5876
$gainers.Cells.Item(62, 1).Value = "🚀"
$gainers.Cells.Item(62, 2).Value = "HITECHGEAR"
$gainers.Cells.Item(62, 3).Value = 3.8587
$gainers.Cells.Item(62, 4).Value = 1.1486
$gainers.Cells.Item(62, 5).Value = 9.9254
$gainers.Cells.Item(63, 1).Value = "🚀"
$gainers.Cells.Item(63, 2).Value = "ORIENTTECH"
$gainers.Cells.Item(63, 3).Value = 3.827
$gainers.Cells.Item(63, 4).Value = 0.5247
$gainers.Cells.Item(63, 5).Value = 32.6784
$gainers.Cells.Item(64, 1).Value = "🚀"
$gainers.Cells.Item(64, 2).Value = "ICRA"
$gainers.Cells.Item(64, 3).Value = 3.7985
$gainers.Cells.Item(64, 4).Value = 4.4793
$gainers.Cells.Item(64, 5).Value = 2.8828
$gainers.Cells.Item(65, 1).Value = "🚀"
$gainers.Cells.Item(65, 2).Value = "SALASAR"
$gainers.Cells.Item(65, 3).Value = 3.7935
$gainers.Cells.Item(65, 4).Value = 4.7872
$gainers.Cells.Item(65, 5).Value = 11.0485
$gainers.Cells.Item(66, 1).Value = "🚀"
$gainers.Cells.Item(66, 2).Value = "NPST"
$gainers.Cells.Item(66, 3).Value = 3.7841
$gainers.Cells.Item(66, 4).Value = -2.0689
$gainers.Cells.Item(66, 5).Value = -3.5677
$gainers.Cells.Item(67, 1).Value = "🚀"
$gainers.Cells.Item(67, 2).Value = "DCW"
$gainers.Cells.Item(67, 3).Value = 3.7544
$gainers.Cells.Item(67, 4).Value = 2.3219
$gainers.Cells.Item(67, 5).Value = -3.9753
$gainers.Cells.Item(68, 1).Value = "🚀"
$gainers.Cells.Item(68, 2).Value = "RHETAN"
$gainers.Cells.Item(68, 3).Value = 3.754
$gainers.Cells.Item(68, 4).Value = 4.178
$gainers.Cells.Item(68, 5).Value = 6.549
$gainers.Cells.Item(69, 1).Value = "🚀"
$gainers.Cells.Item(69, 2).Value = "HINDPETRO"
$gainers.Cells.Item(69, 3).Value = 3.6935
$gainers.Cells.Item(69, 4).Value = 6.9335
$gainers.Cells.Item(69, 5).Value = 5.7397
$gainers.Cells.Item(70, 1).Value = "🚀"
$gainers.Cells.Item(70, 2).Value = "BHARTIHEXA"
$gainers.Cells.Item(70, 3).Value = 3.6718
$gainers.Cells.Item(70, 4).Value = 7.0877
$gainers.Cells.Item(70, 5).Value = 15.3332
$gainers.Cells.Item(71, 1).Value = "🚀"
$gainers.Cells.Item(71, 2).Value = "HLEGLAS"
$gainers.Cells.Item(71, 3).Value = 3.659
$gainers.Cells.Item(71, 4).Value = 8.1155
$gainers.Cells.Item(71, 5).Value = 27.1239
$gainers.Cells.Item(72, 1).Value = "🚀"
$gainers.Cells.Item(72, 2).Value = "RHIM"
$gainers.Cells.Item(72, 3).Value = 3.6544
$gainers.Cells.Item(72, 4).Value = 3.2276
$gainers.Cells.Item(72, 5).Value = 5.1826
$gainers.Cells.Item(73, 1).Value = "🚀"
$gainers.Cells.Item(73, 2).Value = "SHK"
$gainers.Cells.Item(73, 3).Value = 3.6347
$gainers.Cells.Item(73, 4).Value = 2.388
$gainers.Cells.Item(73, 5).Value = -1.932
$gainers.Cells.Item(74, 1).Value = "🚀"
$gainers.Cells.Item(74, 2).Value = "BCLIND"
$gainers.Cells.Item(74, 3).Value = 3.6271
$gainers.Cells.Item(74, 4).Value = 2.2945
$gainers.Cells.Item(74, 5).Value = 0.1728
$gainers.Cells.Item(75, 1).Value = "🚀"
$gainers.Cells.Item(75, 2).Value = "MUKANDLTD"
$gainers.Cells.Item(75, 3).Value = 3.6133
$gainers.Cells.Item(75, 4).Value = 11.9685
$gainers.Cells.Item(75, 5).Value = 9.5508
$gainers.Cells.Item(76, 1).Value = "🚀"
$gainers.Cells.Item(76, 2).Value = "CGPOWER"
$gainers.Cells.Item(76, 3).Value = 3.6125
$gainers.Cells.Item(76, 4).Value = 3.4192
$gainers.Cells.Item(76, 5).Value = 1.0325

# --- Top Losers sheet: refresh leaderboard (drop oldest row, shift up, append new entry) ---
$losers = $wb.Worksheets.Item("Top Losers")
$losers.Rows.Item(70).Delete()
$losers.Cells.Item(70, 1).Value = "📉"
$losers.Cells.Item(70, 2).Value = "JNKINDIA"
$losers.Cells.Item(70, 3).Value = -2.3482
$losers.Cells.Item(70, 4).Value = -2.8371
$losers.Cells.Item(70, 5).Value = 4.2622
$losers.Cells.Item(71, 1).Value = "📉"
$losers.Cells.Item(71, 2).Value = "FCL"
$losers.Cells.Item(71, 3).Value = -2.3453
$losers.Cells.Item(71, 4).Value = -2.616
$losers.Cells.Item(71, 5).Value = -0.02
$losers.Cells.Item(72, 1).Value = "📉"
$losers.Cells.Item(72, 2).Value = "DEEDEV"
$losers.Cells.Item(72, 3).Value = -2.3334
$losers.Cells.Item(72, 4).Value = -6.6528
$losers.Cells.Item(72, 5).Value = -7.4227
$losers.Cells.Item(73, 1).Value = "📉"
$losers.Cells.Item(73, 2).Value = "WEALTH"
$losers.Cells.Item(73, 3).Value = -2.2793
$losers.Cells.Item(73, 4).Value = -3.8356
$losers.Cells.Item(73, 5).Value = -2.7981
$losers.Cells.Item(74, 1).Value = "📉"
$losers.Cells.Item(74, 2).Value = "RATNAMANI"
$losers.Cells.Item(74, 3).Value = -2.2788
$losers.Cells.Item(74, 4).Value = -0.4626
$losers.Cells.Item(74, 5).Value = 0.8712
$losers.Cells.Item(75, 1).Value = "📉"
$losers.Cells.Item(75, 2).Value = "CSBBANK"
$losers.Cells.Item(75, 3).Value = -2.2695
$losers.Cells.Item(75, 4).Value = 2.3137
$losers.Cells.Item(75, 5).Value = 10.6999
$losers.Cells.Item(76, 1).Value = "📉"
$losers.Cells.Item(76, 2).Value = "BBOX"
$losers.Cells.Item(76, 3).Value = -2.2639
$losers.Cells.Item(76, 4).Value = -4.7636
$losers.Cells.Item(76, 5).Value = 5.1528
